$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.282.07"
$ws.Range("E2").Value = "  -2.27%  "

$ws.Range("D3").Value = "3.245.21"
$ws.Range("E3").Value = "  -5.22%  "

$ws.Range("E4").Value = "  +0.40%  "

$ws.Range("D5").Value = "584.72"
$ws.Range("E5").Value = "  -4.10%  "

$ws.Range("D6").Value = "146.85"
$ws.Range("E6").Value = "  -11.54%  "

$ws.Range("D8").Value = "3.240.96"
$ws.Range("E8").Value = "  -5.27%  "

$ws.Range("D9").Value = "0.536"
$ws.Range("E9").Value = "  -9.28%  "

$ws.Range("D10").Value = "0.167"
$ws.Range("E10").Value = "  -12.91%  "

$ws.Range("D11").Value = "6.64"
$ws.Range("E11").Value = "  -4.29%  "

$ws.Range("D12").Value = "0.497"
$ws.Range("E12").Value = "  -11.10%  "

$ws.Range("D13").Value = "0.0000242"
$ws.Range("E13").Value = "  -9.27%  "

$ws.Range("D14").Value = "37.61"
$ws.Range("E14").Value = "  -13.93%  "

$ws.Range("D15").Value = "3.778.00"
$ws.Range("E15").Value = "  -4.93%  "

$ws.Range("D16").Value = "67.361.49"
$ws.Range("E16").Value = "  -2.20%  "

$ws.Range("D17").Value = "3.259.73"
$ws.Range("E17").Value = "  -4.71%  "

$ws.Range("E18").Value = "  -6.01%  "

$ws.Range("D19").Value = "519.13"
$ws.Range("E19").Value = "  -9.95%  "

$ws.Range("D20").Value = "6.98"
$ws.Range("E20").Value = "  -13.39%  "

$ws.Range("D21").Value = "14.70"
$ws.Range("E21").Value = "  -13.53%  "

$ws.Range("D22").Value = "0.742"
$ws.Range("E22").Value = "  -11.53%  "

$ws.Range("D23").Value = "7.67"
$ws.Range("E23").Value = "  -13.06%  "

$ws.Range("D24").Value = "84.60"

$ws.Range("D25").Value = "13.20"
$ws.Range("E25").Value = "  -11.82%  "

$ws.Range("D26").Value = "0.998"
$ws.Range("E26").Value = "  -0.28%  "

$ws.Range("D27").Value = "3.18"
$ws.Range("E27").Value = "  -11.59%  "

$ws.Range("D28").Value = "2.11"
$ws.Range("E28").Value = "  -11.68%  "

$ws.Range("D29").Value = "7.85"
$ws.Range("E29").Value = "  -7.50%  "

$ws.Range("D30").Value = "28.52"
$ws.Range("E30").Value = "  -12.16%  "

$ws.Range("D31").Value = "1.17"
$ws.Range("E31").Value = "  -4.47%  "

$ws.Range("D32").Value = "2.61"
$ws.Range("E32").Value = "  -4.90%  "

$ws.Range("D33").Value = "6.43"
$ws.Range("E33").Value = "  -16.92%  "

$ws.Range("B34").Value = "FirstDigitalUSD"
$ws.Range("C34").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D34").Value = "1.01"
$ws.Range("E34").Value = "  +0.37%  "

$ws.Range("D35").Value = "5.59"
$ws.Range("E35").Value = "  -14.02%  "

$ws.Range("B36").Value = "OKB"
$ws.Range("C36").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D36").Value = "55.51"
$ws.Range("E36").Value = "  -0.38%  "

$ws.Range("D37").Value = "502.05"
$ws.Range("E37").Value = "  -13.51%  "

$ws.Range("D38").Value = "0.0434"
$ws.Range("E38").Value = "  -6.36%  "

$ws.Range("D39").Value = "0.0838"
$ws.Range("E39").Value = "  -11.42%  "

$ws.Range("E40").Value = "  -11.44%  "

$ws.Range("E41").Value = "  -15.55%  "

$ws.Range("D42").Value = "2.895.38"
$ws.Range("E42").Value = "  -10.08%  "

$ws.Range("D43").Value = "2.66"
$ws.Range("E43").Value = "  -13.20%  "

$ws.Range("D44").Value = "0.261"
$ws.Range("E44").Value = "  -10.32%  "

$ws.Range("D45").Value = "2.16"
$ws.Range("E45").Value = "  -8.59%  "

$ws.Range("B46").Value = "USDe"
$ws.Range("C46").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D46").Value = "1.00"
$ws.Range("E46").Value = "  -0.11%  "

$ws.Range("D47").Value = "26.10"
$ws.Range("E47").Value = "  -15.20%  "

$ws.Range("B48").Value = "PEPE"
$ws.Range("C48").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D48").Value = "0.0₃0556"
$ws.Range("E48").Value = "  -17.20%  "

$ws.Range("D49").Value = "123.76"
$ws.Range("E49").Value = "  -6.21%  "

$ws.Range("B50").Value = "Stellar"
$ws.Range("C50").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D50").Value = "0.112"
$ws.Range("E50").Value = "  -10.53%  "

$ws.Range("B51").Value = "ThetaToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D51").Value = "2.26"
$ws.Range("E51").Value = "  -17.82%  "
